$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their text representation (some values
# like "1.006" would otherwise be auto-converted to numbers by Excel).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.409.36'
$ws.Range("E2").Value = '  -0.02%  '
$ws.Range("D3").Value = '1.820.84'
$ws.Range("E3").Value = '  -0.39%  '
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.42%  '
$ws.Range("D5").Value = '316.38'
$ws.Range("E5").Value = '  +0.49%  '
$ws.Range("E6").Value = '  +0.36%  '
$ws.Range("D7").Value = '0.5131'
$ws.Range("E7").Value = '  -0.58%  '
$ws.Range("D8").Value = '0.3847'
$ws.Range("E8").Value = '  -1.74%  '
$ws.Range("D9").Value = '0.08219'
$ws.Range("E9").Value = '  +7.26%  '
$ws.Range("D10").Value = '1.117'
$ws.Range("E10").Value = '  +0.76%  '
$ws.Range("D11").Value = '41.89'
$ws.Range("E11").Value = '  +0.05%  '
$ws.Range("D12").Value = '6.352'
$ws.Range("E12").Value = '  +1.15%  '
$ws.Range("D13").Value = '21.05'
$ws.Range("D14").Value = '1.006'
$ws.Range("E14").Value = '  +0.41%  '
$ws.Range("D15").Value = '7.458'
$ws.Range("E15").Value = '  -1.18%  '
$ws.Range("D16").Value = '1.821.05'
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("D17").Value = '93.94'
$ws.Range("E17").Value = '  +0.25%  '
$ws.Range("D18").Value = '0.00001112'
$ws.Range("E18").Value = '  +1.97%  '
$ws.Range("D19").Value = '0.06618'
$ws.Range("E19").Value = '  -1.53%  '
$ws.Range("D20").Value = '17.74'
$ws.Range("E20").Value = '  +0.39%  '
$ws.Range("D21").Value = '1.005'
$ws.Range("E21").Value = '  +0.42%  '
$ws.Range("D22").Value = '6.040'
$ws.Range("E22").Value = '  -2.28%  '
$ws.Range("D23").Value = '28.457.61'
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("D24").Value = '11.46'
$ws.Range("E24").Value = '  +2.71%  '
$ws.Range("D25").Value = '2.246'
$ws.Range("E25").Value = '  -0.38%  '
$ws.Range("D26").Value = '160.19'
$ws.Range("E26").Value = '  +2.10%  '
$ws.Range("D27").Value = '20.96'
$ws.Range("E27").Value = '  +1.66%  '
$ws.Range("D28").Value = '2.027.97'
$ws.Range("E28").Value = '  -0.49%  '
$ws.Range("D29").Value = '2.396'
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("D30").Value = '125.41'
$ws.Range("E30").Value = '  +0.40%  '
$ws.Range("D31").Value = '0.1096'
$ws.Range("E31").Value = '  +0.68%  '
$ws.Range("D32").Value = '1.088'
$ws.Range("E32").Value = '  -2.51%  '
$ws.Range("D33").Value = '5.705'
$ws.Range("E33").Value = '  +0.53%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '3.686'
$ws.Range("E34").Value = '  +0.66%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.07428'
$ws.Range("E35").Value = '  +6.00%  '
$ws.Range("D36").Value = '12.45'
$ws.Range("E36").Value = '  +10.76%  '
$ws.Range("D37").Value = '0.2213'
$ws.Range("E37").Value = '  -0.40%  '
$ws.Range("D38").Value = '0.02347'
$ws.Range("E38").Value = '  +1.13%  '
$ws.Range("D39").Value = '5.194'
$ws.Range("E39").Value = '  +0.73%  '
$ws.Range("D40").Value = '8.812'
$ws.Range("E40").Value = '  -1.13%  '
$ws.Range("D41").Value = '0.6338'
$ws.Range("E41").Value = '  +0.84%  '
$ws.Range("D42").Value = '1.181'
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("D43").Value = '1.389'
$ws.Range("E43").Value = '  -0.19%  '
$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").Value = '0.6169'
$ws.Range("E44").Value = '  +4.56%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '13.51'
$ws.Range("E45").Value = '  +0.79%  '
$ws.Range("D46").Value = '3.817'
$ws.Range("E46").Value = '  +2.89%  '
$ws.Range("D47").Value = '127.34'
$ws.Range("E47").Value = '  +2.16%  '
$ws.Range("D48").Value = '1.996'
$ws.Range("E48").Value = '  +1.01%  '
$ws.Range("D49").Value = '1.203'
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("D50").Value = '0.06915'
$ws.Range("E50").Value = '  -0.29%  '
$ws.Range("D51").Value = '1.073'
$ws.Range("E51").Value = '  +1.05%  '
